# Update the "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" sheets to reflect the freshly regenerated data snapshot.
# Each listed row's column F value is incremented to match the new output.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 43
    8  = 119
    13 = 1099
    20 = 55
    26 = 60
    31 = 3923
    38 = 21
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
